$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 5
$ws.Range("P2").Value = 35
$ws.Range("R2").Value = 20
$ws.Range("U2").Value = 15
$ws.Range("V2").Value = 15
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 6
$ws.Range("R3").Value = 40
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 7
$ws.Range("O4").Value = 35
$ws.Range("P4").Value = 35
$ws.Range("Q4").Value = 25
$ws.Range("S4").Value = 15
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 1
$ws.Range("N5").Value = 30
$ws.Range("Q5").Value = 35
$ws.Range("F6").Value = 2
$ws.Range("J6").Value = 2
$ws.Range("M6").Value = 65
$ws.Range("N6").Value = 20
$ws.Range("P6").Value = 45
$ws.Range("Q6").Value = 45
$ws.Range("R6").Value = 10
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 9
$ws.Range("G7").Value = 5
$ws.Range("J7").Value = 2
$ws.Range("O7").Value = 45
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 6
$ws.Range("J8").Value = 4
$ws.Range("S8").Value = 20
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 5
$ws.Range("I9").Value = 2
$ws.Range("N9").Value = 50
$ws.Range("P9").Value = 40
$ws.Range("R9").Value = 10
$ws.Range("S9").Value = 35
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 5
$ws.Range("M10").Value = 55
$ws.Range("N10").Value = 40
$ws.Range("O10").Value = 30
$ws.Range("P10").Value = 30
$ws.Range("R10").Value = 5
$ws.Range("S10").Value = 45
$ws.Range("D11").Value = 17
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 4
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 3
$ws.Range("M11").Value = 60
$ws.Range("O11").Value = 40
$ws.Range("R11").Value = 20
$ws.Range("D12").Value = 16
$ws.Range("F12").Value = 7
$ws.Range("G12").Value = 9
$ws.Range("J12").Value = 6
$ws.Range("N12").Value = 30
$ws.Range("O12").Value = 35
$ws.Range("R12").Value = 10
$ws.Range("D13").Value = 18
$ws.Range("I13").Value = 7
$ws.Range("N13").Value = 45
$ws.Range("O13").Value = 35
$ws.Range("R13").Value = 35
$ws.Range("S13").Value = 15
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 5
$ws.Range("I14").Value = 4
$ws.Range("R14").Value = 20
$ws.Range("S14").Value = 10
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 6
$ws.Range("M15").Value = 90
$ws.Range("O15").Value = 25
$ws.Range("P15").Value = 40
$ws.Range("S15").Value = 10

$ws.Range("E9").Select()
